$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D keep their exact text
# representation (including trailing/leading zeros) instead of being
# auto-converted to floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.169.20'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.977.42'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '565.66'
$ws.Range("E5").Value = '  +1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.81'
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.973.75'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("E11").Value = '  +8.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.49'
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.469.65'
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.06'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.976.71'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '59.161.39'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '433.73'
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.60'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.723'
$ws.Range("E22").Value = '  +1.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.02'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.97'
$ws.Range("E24").Value = '  -3.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.90'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.21'
$ws.Range("E27").Value = '  +5.60%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.69'
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.63'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("E33").Value = '  +4.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.985'
$ws.Range("E34").Value = '  +2.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.86'
$ws.Range("E35").Value = '  +1.60%  '
$ws.Range("E36").Value = '  +2.72%  '
$ws.Range("E37").Value = '  -2.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.47'
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.72'
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '393.90'
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0350'
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.702.32'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.104'
$ws.Range("E44").Value = '  -3.49%  '
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.56'
$ws.Range("E47").Value = '  -2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.29'
$ws.Range("E48").Value = '  +10.72%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.96'
$ws.Range("E50").Value = '  -2.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.04'
$ws.Range("E51").Value = '  -1.05%  '
